$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Table 1 - Bill header (Date / ID No / Name / TP No)
# ---------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$t1.Cell(1, 1).Range.Text = "දිනය"
$t1.Cell(1, 2).Range.Text = "2020-12-28 13:07:37"
$t1.Cell(2, 1).Range.Text = "ජා. අංකය"
$t1.Cell(2, 2).Range.Text = "123"
$t1.Cell(3, 1).Range.Text = "නම"
$t1.Cell(3, 2).Range.Text = "Test"
$t1.Cell(4, 1).Range.Text = "දු. අංකය"
$t1.Cell(4, 2).Range.Text = "123"

# ---------------------------------------------------------------
# Table 2 - Rented Item details
# ---------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$t2.Cell(1, 1).Range.Text = "භාණ්ඩ වර්ගය"
$t2.Cell(1, 2).Range.Text = "කුලියට ගත් දිනය"
$t2.Cell(1, 3).Range.Text = "දින ගණන"
$t2.Cell(1, 4).Range.Text = "ප්‍රමාණය"
$t2.Cell(1, 5).Range.Text = "දිනකට කුලිය"
$t2.Cell(1, 6).Range.Text = "මුදල"

# Remove the only item row (Poker / 2020-12-10 / 11.0 / 2 / 2,000.00 / 44,000.00)
$t2.Rows.Item(2).Delete()

# Re-fetch the table/row after the structural edit and update the Total amount
$t2 = $d.Tables.Item(2)
$t2.Cell(2, 6).Range.Text = "      0.00"

# ---------------------------------------------------------------
# Table 3 - Payment Details
# ---------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$t3.Cell(1, 1).Range.Text = "හිඟ මුදල"
$t3.Cell(1, 2).Range.Text = "Rs.  13,000.00"
$t3.Cell(2, 1).Range.Text = "මෙම බිල්පතෙහි වටිනාකම"
$t3.Cell(2, 2).Range.Text = "Rs.       0.00"
$t3.Cell(3, 1).Range.Text = "මුළු මුදල"
$t3.Cell(3, 2).Range.Text = "            Rs.  13,000.00"
$t3.Cell(4, 1).Range.Text = "ගෙවීම්"
$t3.Cell(4, 2).Range.Text = "-(Rs.  13,000.00)"
$t3.Cell(5, 1).Range.Text = "ගෙවිය යුතු වටිනාකම"
$t3.Cell(5, 2).Range.Text = "            Rs.       0.00"

Write-Output "done"
